$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows with text values, preserving leading zeros.
# Leading apostrophe forces text entry (like Excel's own text-number
# handling); ClearFormats afterwards drops the quote-prefix formatting
# that gets attached so the cell keeps the workbook's default style,
# matching the other plain data cells in the column.
$ws.Range("A5").Value = "'011748367"
$ws.Range("A5").ClearFormats()

$ws.Range("A6").Value = "'011748625"
$ws.Range("A6").ClearFormats()
